$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2206.8572
$ws.Range("I31").Value = 74.666664
$ws.Range("K31").Value = 223.999992
$ws.Range("M31").Value = 6.000008000000008
$ws.Range("H32").Value = 800
$ws.Range("J32").Value = 800
$ws.Range("L32").Value = 800
$ws.Range("N32").Value = -1452
$ws.Range("H41").Value = 2487
$ws.Range("I41").Value = 1999.5
$ws.Range("J41").Value = 2974.5
$ws.Range("K41").Value = 1999.5
$ws.Range("L41").Value = 2974.5
$ws.Range("M41").Value = -1559.5
$ws.Range("N41").Value = -3854.5
$ws.Range("H64").Value = 7022.727
$ws.Range("J64").Value = 7466.5
$ws.Range("L64").Value = 7466.5
$ws.Range("N64").Value = -7962.5
$ws.Range("H67").Value = 7022.727
$ws.Range("J67").Value = 7466.5
$ws.Range("L67").Value = 7466.5
$ws.Range("N67").Value = -9182.5
$ws.Range("H74").Value = 16200
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -64
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 16200
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -320
$ws.Range("N77").Value = -109360
$ws.Range("H80").Value = 456.28125
$ws.Range("I80").Value = 272.35715
$ws.Range("J80").Value = 599.3333
$ws.Range("K80").Value = 817.0714499999999
$ws.Range("L80").Value = 1797.9999
$ws.Range("M80").Value = 180.9285500000001
$ws.Range("N80").Value = -3793.9999
$ws.Range("H83").Value = 456.28125
$ws.Range("I83").Value = 272.35715
$ws.Range("J83").Value = 599.3333
$ws.Range("K83").Value = 2451.21435
$ws.Range("L83").Value = 5393.9997
$ws.Range("M83").Value = 2540.78565
$ws.Range("N83").Value = -15377.9997
$ws.Range("H133").Value = 140780
$ws.Range("J133").Value = 140780
$ws.Range("L133").Value = 140780
$ws.Range("N133").Value = -150900
$ws.Range("H137").Value = 2429.5
$ws.Range("I137").Value = 1358
$ws.Range("K137").Value = 4074
$ws.Range("M137").Value = -1524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2863981.2
$ws.Range("I32").Value = 3687.3333
$ws.Range("J32").Value = 12517473
$ws.Range("K32").Value = 3687.3333
$ws.Range("L32").Value = 12517473
$ws.Range("M32").Value = -3400.3333
$ws.Range("N32").Value = -12518047
$ws.Range("H45").Value = 3032.25
$ws.Range("J45").Value = 5666.6665
$ws.Range("L45").Value = 5666.6665
$ws.Range("N45").Value = -6420.6665
$ws.Range("H61").Value = 3177.7778
$ws.Range("I61").Value = 3177.7778
$ws.Range("K61").Value = 3177.7778
$ws.Range("M61").Value = -2965.7778
$ws.Range("H74").Value = 3540.24
$ws.Range("I74").Value = 2926.3333
$ws.Range("J74").Value = 6763.25
$ws.Range("K74").Value = 2926.3333
$ws.Range("L74").Value = 6763.25
$ws.Range("M74").Value = -2052.3333
$ws.Range("N74").Value = -8511.25
$ws.Range("H77").Value = 3540.24
$ws.Range("I77").Value = 2926.3333
$ws.Range("J77").Value = 6763.25
$ws.Range("K77").Value = 14631.6665
$ws.Range("L77").Value = 33816.25
$ws.Range("M77").Value = -10263.6665
$ws.Range("N77").Value = -42552.25
$ws.Range("H110").Value = 2632.75
$ws.Range("I110").Value = 2478
$ws.Range("K110").Value = 2478
$ws.Range("M110").Value = -433
$ws.Range("H125").Value = 55951.5
$ws.Range("J125").Value = 55951.5
$ws.Range("L125").Value = 55951.5
$ws.Range("N125").Value = -65791.5
$ws.Range("H132").Value = 6199.8335
$ws.Range("I132").Value = 2606
$ws.Range("K132").Value = 7818
$ws.Range("M132").Value = -5288
$ws.Range("H136").Value = 3177.7778
$ws.Range("I136").Value = 3177.7778
$ws.Range("K136").Value = 9533.3334
$ws.Range("M136").Value = -6983.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2238.5
$ws.Range("I20").Value = 2442.5881
$ws.Range("K20").Value = 2442.5881
$ws.Range("M20").Value = -2195.5881

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 204.83333
$ws.Range("I10").Value = 184.2
$ws.Range("K10").Value = 184.2
$ws.Range("M10").Value = -45.19999999999999
$ws.Range("H58").Value = 4183.385
$ws.Range("I58").Value = 885
$ws.Range("J58").Value = 6244.875
$ws.Range("K58").Value = 885
$ws.Range("L58").Value = 6244.875
$ws.Range("M58").Value = -682
$ws.Range("N58").Value = -6650.875
$ws.Range("H60").Value = 93
$ws.Range("I60").Value = 93
$ws.Range("K60").Value = 93
$ws.Range("M60").Value = 418
$ws.Range("H132").Value = 4354.8
$ws.Range("J132").Value = 5948.5557
$ws.Range("L132").Value = 17845.6671
$ws.Range("N132").Value = -22905.6671
$ws.Range("H136").Value = 4183.385
$ws.Range("I136").Value = 885
$ws.Range("J136").Value = 6244.875
$ws.Range("K136").Value = 2655
$ws.Range("L136").Value = 18734.625
$ws.Range("M136").Value = -105
$ws.Range("N136").Value = -23834.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 28.619047
$ws.Range("I2").Value = 28.75
$ws.Range("K2").Value = 172.5
$ws.Range("M2").Value = -59.5
$ws.Range("H12").Value = 377.81818
$ws.Range("I12").Value = 191.66667
$ws.Range("J12").Value = 447.625
$ws.Range("K12").Value = 575.00001
$ws.Range("L12").Value = 1342.875
$ws.Range("M12").Value = -402.00001
$ws.Range("N12").Value = -1688.875
$ws.Range("H38").Value = 311.90475
$ws.Range("I38").Value = 326.25
$ws.Range("K38").Value = 978.75
$ws.Range("M38").Value = -631.75
$ws.Range("H60").Value = 990
$ws.Range("I60").Value = 299.6875
$ws.Range("J60").Value = 2567.8572
$ws.Range("K60").Value = 899.0625
$ws.Range("L60").Value = 7703.571599999999
$ws.Range("M60").Value = -648.0625
$ws.Range("N60").Value = -8205.571599999999
$ws.Range("H98").Value = 253.16667
$ws.Range("I98").Value = 131
$ws.Range("J98").Value = 497.5
$ws.Range("K98").Value = 393
$ws.Range("L98").Value = 1492.5
$ws.Range("M98").Value = 1105
$ws.Range("N98").Value = -4488.5
$ws.Range("H129").Value = 1970.091
$ws.Range("J129").Value = 7516
$ws.Range("L129").Value = 22548
$ws.Range("N129").Value = -32548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H55").Value = 4802.25
$ws.Range("I55").Value = 4971.5
$ws.Range("J55").Value = 4294.5
$ws.Range("K55").Value = 4971.5
$ws.Range("L55").Value = 4294.5
$ws.Range("M55").Value = -4644.5
$ws.Range("N55").Value = -4948.5
$ws.Range("H102").Value = 1083.45
$ws.Range("I102").Value = 1103.1578
$ws.Range("K102").Value = 1103.1578
$ws.Range("M102").Value = 518.8422
$ws.Range("H122").Value = 455981.53
$ws.Range("I122").Value = 501429.7
$ws.Range("K122").Value = 1504289.1
$ws.Range("M122").Value = -1501839.1
$ws.Range("H126").Value = 5499.6665
$ws.Range("I126").Value = 5499
$ws.Range("K126").Value = 16497
$ws.Range("M126").Value = -14027

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8064.6665
$ws.Range("I7").Value = 5832.6665
$ws.Range("J7").Value = 8808.666999999999
$ws.Range("K7").Value = 5832.6665
$ws.Range("L7").Value = 8808.666999999999
$ws.Range("M7").Value = -5720.6665
$ws.Range("N7").Value = -9032.666999999999
$ws.Range("H22").Value = 940
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -1190
$ws.Range("H27").Value = 940
$ws.Range("J27").Value = 600
$ws.Range("L27").Value = 600
$ws.Range("N27").Value = -814
$ws.Range("H40").Value = 5178.6924
$ws.Range("I40").Value = 3176.111
$ws.Range("K40").Value = 3176.111
$ws.Range("M40").Value = -3040.111
$ws.Range("H126").Value = 8064.6665
$ws.Range("I126").Value = 5832.6665
$ws.Range("J126").Value = 8808.666999999999
$ws.Range("K126").Value = 17497.9995
$ws.Range("L126").Value = 26426.001
$ws.Range("M126").Value = -15027.9995
$ws.Range("N126").Value = -31366.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2167
$ws.Range("J7").Value = 2250
$ws.Range("L7").Value = 2250
$ws.Range("N7").Value = -2476
$ws.Range("H102").Value = 337000
$ws.Range("J102").Value = 337000
$ws.Range("L102").Value = 337000
$ws.Range("N102").Value = -343490
$ws.Range("H113").Value = 598.5
$ws.Range("I113").Value = 711.7273
$ws.Range("J113").Value = 349.4
$ws.Range("K113").Value = 2135.1819
$ws.Range("L113").Value = 1048.2
$ws.Range("M113").Value = 34.81809999999996
$ws.Range("N113").Value = -5388.2
$ws.Range("H122").Value = 2184.25
$ws.Range("I122").Value = 2139.1428
$ws.Range("K122").Value = 6417.428400000001
$ws.Range("M122").Value = -3967.428400000001
$ws.Range("H126").Value = 4293.0835
$ws.Range("I126").Value = 1669.6666
$ws.Range("J126").Value = 6916.5
$ws.Range("K126").Value = 5008.9998
$ws.Range("L126").Value = 20749.5
$ws.Range("M126").Value = -2538.9998
$ws.Range("N126").Value = -25689.5
